# Auto-generated edit script to update cryptos list data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "90.224.95"
$ws.Range("E2").Value = "  -0.51%  "

# Row 3
$ws.Range("D3").Value = "3.056.60"
$ws.Range("E3").Value = "  -1.67%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.93"
$ws.Range("E5").Value = "  +2.80%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "616.47"
$ws.Range("E6").Value = "  -2.16%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.12"
$ws.Range("E7").Value = "  +7.15%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.363"
$ws.Range("E8").Value = "  +0.51%  "

# Row 9
$ws.Range("E9").Value = "  +0.05%  "

# Row 10
$ws.Range("D10").Value = "3.059.85"
$ws.Range("E10").Value = "  -1.57%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.734"
$ws.Range("E11").Value = "  +2.91%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.200"
$ws.Range("E12").Value = "  +2.86%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("E13").Value = "  +0.57%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.78"
$ws.Range("E14").Value = "  -4.59%  "

# Row 15
$ws.Range("D15").Value = "90.793.43"
$ws.Range("E15").Value = "  +0.22%  "

# Row 17
$ws.Range("D17").Value = "3.647.14"
$ws.Range("E17").Value = "  -1.08%  "

# Row 18
$ws.Range("D18").Value = "3.099.59"
$ws.Range("E18").Value = "  -1.96%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.67"
$ws.Range("E19").Value = "  -1.72%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.33"
$ws.Range("E20").Value = "  +1.79%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000210"
$ws.Range("E21").Value = "  +1.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.73"
$ws.Range("E22").Value = "  +3.73%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "438.45"
$ws.Range("E23").Value = "  -0.49%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.00"
$ws.Range("E24").Value = "  +1.01%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "90.38"
$ws.Range("E25").Value = "  +3.17%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.56"
$ws.Range("E26").Value = "  -5.28%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.71"
$ws.Range("E27").Value = "  -5.20%  "

# Row 28
$ws.Range("E28").Value = "  -0.95%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.06%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.245"
$ws.Range("E30").Value = "  +26.94%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.182"
$ws.Range("E31").Value = "  +14.36%  "

# Row 32
$ws.Range("E32").Value = "  +13.83%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.09"
$ws.Range("E33").Value = "  -3.69%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.111"
$ws.Range("E34").Value = "  +31.92%  "

# Row 35
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.165"
$ws.Range("E35").Value = "  +11.47%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.67"
$ws.Range("E36").Value = "  +8.81%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.22"
$ws.Range("E37").Value = "  -0.05%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.20"
$ws.Range("E38").Value = "  +29.26%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.90"
$ws.Range("E39").Value = "  -0.61%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "489.87"
$ws.Range("E40").Value = "  -3.47%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.58"
$ws.Range("E41").Value = "  -5.04%  "

# Row 42
$ws.Range("E42").Value = "  +0.35%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.414"
$ws.Range("E43").Value = "  +1.01%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.09"
$ws.Range("E44").Value = "  -0.35%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "153.94"
$ws.Range("E46").Value = "  +1.87%  "

# Row 47
$ws.Range("E47").Value = "  -1.15%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.678"
$ws.Range("E48").Value = "  -0.73%  "

# Row 49
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.03"
$ws.Range("E49").Value = "  -2.19%  "

# Row 50
$ws.Range("E50").Value = "  -1.16%  "

# Row 51
$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.41"
$ws.Range("E51").Value = "  -1.29%  "

